{"js": "// Update the date heading paragraph.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst firstPara = paras.items[0];\nfirstPara.load(\"text\");\nawait context.sync();\nif (firstPara.text.trim() === \"2024-11-27 Wednesday\") {\n  firstPara.insertText(\"2024-11-28 Thursday\", \"Replace\");\n}\n\n// Update the multiplication-problem table. Each \"problem\" row is followed\n// by four blank rows (for the student's work); only the problem rows carry\n// text, five two-digit-by-two-digit problems per row.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = {\n  0: [\"85\u00d773=\", \"44\u00d725=\", \"98\u00d771=\", \"16\u00d791=\", \"69\u00d765=\"],\n  4: [\"82\u00d768=\", \"30\u00d771=\", \"47\u00d749=\", \"87\u00d714=\", \"67\u00d753=\"],\n  9: [\"31\u00d743=\", \"90\u00d761=\", \"24\u00d713=\", \"19\u00d737=\", \"62\u00d774=\"],\n  14: [\"58\u00d719=\", \"58\u00d767=\", \"12\u00d759=\", \"58\u00d727=\", \"16\u00d737=\"],\n  19: [\"90\u00d730=\", \"71\u00d770=\", \"81\u00d779=\", \"23\u00d792=\", \"55\u00d760=\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const r = parseInt(rowIndex, 10);\n  const rowVals = newValues[r];\n  for (let c = 0; c < rowVals.length; c++) {\n    table.getCell(r, c).value = rowVals[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading paragraph.\n$p = $d.Paragraphs.Item(1)\nif ($p.Range.Text.Trim() -eq \"2024-11-27 Wednesday\") {\n    $p.Range.Text = \"2024-11-28 Thursday\"\n}\n\n# Update the multiplication-problem table. Each \"problem\" row is followed by\n# four blank rows (for the student's work); only the problem rows carry\n# text, five two-digit-by-two-digit problems per row. Row/column numbers\n# below are 1-based, matching Word's Table.Cell(row, column) addressing.\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"85\u00d773=\", \"44\u00d725=\", \"98\u00d771=\", \"16\u00d791=\", \"69\u00d765=\")\n    5  = @(\"82\u00d768=\", \"30\u00d771=\", \"47\u00d749=\", \"87\u00d714=\", \"67\u00d753=\")\n    10 = @(\"31\u00d743=\", \"90\u00d761=\", \"24\u00d713=\", \"19\u00d737=\", \"62\u00d774=\")\n    15 = @(\"58\u00d719=\", \"58\u00d767=\", \"12\u00d759=\", \"58\u00d727=\", \"16\u00d737=\")\n    20 = @(\"90\u00d730=\", \"71\u00d770=\", \"81\u00d779=\", \"23\u00d792=\", \"55\u00d760=\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowVals = $newValues[$rowIndex]\n    for ($c = 1; $c -le $rowVals.Length; $c++) {\n        $t.Cell($rowIndex, $c).Range.Text = $rowVals[$c - 1]\n    }\n}\n"}
